$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# The run "for c in data%}" is split into three runs: "for c in ", "data"
# and "%}", with a collapsed _GoBack bookmark inserted right after "data"
# (i.e. between "data" and "%}"), mimicking Word re-typing the word
# "data" at that location.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("data", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found1) {
    # Mark the start of "data" with a temporary collapsed bookmark so the
    # preceding run ("for c in ") is split off from "data".
    $startPoint = $d.Range($rng1.Start, $rng1.Start)
    $d.Bookmarks.Add("TempSplitMark", $startPoint)

    # Mark the real _GoBack bookmark (collapsed) right after "data", which
    # splits "data" from the trailing "%}" run.
    $endPoint = $d.Range($rng1.End, $rng1.End)
    $d.Bookmarks.Add("_GoBack", $endPoint)

    # Remove the temporary bookmark - the run split it caused remains.
    $d.Bookmarks("TempSplitMark").Delete()
}

# --- Change 2 ---------------------------------------------------------
# The table header cell text was previously split into "P" + bookmark +
# "rice"; it is now a single run "Price" with the bookmark removed.
$rng2 = $d.Content
$null = $rng2.Find.Execute("Price", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "Price", 2)
